$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredential")

# Update the data table (B2:E5) with new values
$ws.Range("B2").Value = "GoogleChrome"
$ws.Range("C2").Value = "mngrINVALID"
$ws.Range("D2").Value = "udAhydy"
$ws.Range("E2").Value = "FAIL"

$ws.Range("B3").Value = "MozilaFireFox"
$ws.Range("C3").Value = "mngr243120"
$ws.Range("D3").Value = "udAhydy"
$ws.Range("E3").Value = "PASS"

$ws.Range("B4").Value = "MozilaFireFox"
$ws.Range("C4").Value = "mng243120"
$ws.Range("D4").Value = "udAhydy"
$ws.Range("E4").Value = "PASS"

$ws.Range("B5").Value = "InternetExplore"
$ws.Range("C5").Value = "INVALID"
$ws.Range("D5").Value = "INVALIDJG44"
$ws.Range("E5").Value = "FAIL"

# Update column B width (results in stored width="14" in the OOXML)
$ws.Columns.Item(2).ColumnWidth = 13.15

# Update the selected cell
$ws.Range("G12").Select()
